$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293, shifting existing rows 293-364 down to 294-365.
$ws.Rows.Item(293).Insert()

# Populate the new row 293 with its data.
$ws.Cells.Item(293, 1).Value = 9
$ws.Cells.Item(293, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(293, 3).Value = "Metropolitana"
$ws.Cells.Item(293, 4).Value = "2023-01-06"
$ws.Cells.Item(293, 5).Value = 13
$ws.Cells.Item(293, 6).Value = 300000001
$ws.Cells.Item(293, 7).Value = "Rabanito"
$ws.Cells.Item(293, 8).Value = "Sin especificar"
$ws.Cells.Item(293, 9).Value = "Primera"
$ws.Cells.Item(293, 10).Value = 7000
$ws.Cells.Item(293, 11).Value = 3000
$ws.Cells.Item(293, 12).Value = 3000
$ws.Cells.Item(293, 13).Value = 3000
$ws.Cells.Item(293, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(293, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(293, 16).Value = 30
$ws.Cells.Item(293, 17).Value = 100
$ws.Cells.Item(293, 18).Value = "Hortaliza"
